# Add produce Logic of SLG Building
# - Remove the "ID" row (row 2: Id=ID, Type=string, Desc=建筑ID)
# - Change remaining Prefab/NormalStateFunc/UpStateFunc/Desc rows' Type from int -> string
# - Data validation list range shrinks from F9:F1048576 to F8:F1048576, prompts re-enabled
# - Selection / absolute path bookkeeping

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the obsolete "ID" row (current row 2), shifting rows 3-8 up to 2-7
$ws.Rows.Item(2).Delete()

# The Prefab / NormalStateFunc / UpStateFunc / Desc rows are now rows 4-7;
# flip their "Type" column (B) from int -> string
$ws.Range("B4:B7").Value = "string"

# Data validation: the TRUE/FALSE list range shifted with the row delete
# (it now covers F8:F1048575); drop it and re-add it anchored at F8:F1048576
# with prompts re-enabled (matches the target sqref + dropped disablePrompts)
$ws.Range("F8:F1048575").Validation.Delete()
$ws.Range("F8:F1048576").Validation.Add(3, 1, 1, """TRUE,FALSE""")
$ws.Range("F8:F1048576").Validation.IgnoreBlank = $true
$ws.Range("F8:F1048576").Validation.InCellDropdown = $true
$ws.Range("F8:F1048576").Validation.ShowInput = $true
$ws.Range("F8:F1048576").Validation.ShowError = $true

# Restore the cursor position that was captured on save
$ws.Range("G14").Select()
